$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("Q1").Value = "Código lote"
$ws.Range("R1").Value = "Fec. Vencimiento"

# Row 3 edits: product name + code changes
$ws.Range("A3").Value = "Pepsi Coca cola"
$ws.Range("B3").Value = "P1212"

# New data for row 3
$ws.Range("Q3").Value = 11223344
$ws.Range("R3").Value = 44123
$ws.Range("R3").NumberFormat = "mm-dd-yy"

# Selection change
$ws.Range("D5").Select()

$ws.Columns.Item(18).ColumnWidth = 15.67
